# Task 2 - identify the column positions of the credential headers
# ("Account No" and "PAN Number") in the dataset, then move the PAN
# Number column so the layout becomes: Name | Account No | Subject |
# Marks | PAN Number  (i.e. the PAN column is rotated to the end,
# single step, using the detected column positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- detect the header columns dynamically ---------------------------
$acctCol = 0
$panCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($header -eq "Account No") {
        $acctCol = $c
    }
    if ($header -eq "PAN Number") {
        $panCol = $c
    }
}

Write-Host "Detected 'Account No' column:" $acctCol
Write-Host "Detected 'PAN Number' column:" $panCol

# --- randomise / relocate the credentials column in a single step ----
# Move the PAN Number column to sit immediately after the last column,
# which shifts every column in between one position to the left
# (Account No, Subject, Marks all move left by one) and the PAN Number
# column becomes the new last column.
if ($panCol -gt 0 -and $panCol -lt $lastCol) {
    $ws.Columns.Item($panCol).Cut()
    $ws.Columns.Item($lastCol + 1).Insert()
}

# --- restore the active selection recorded for this sheet ------------
$ws.Range("H13").Select()
